$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 31 (2976 / 대림강변타운), shifting
# that row and everything below it down by one. Fill the new row with the
# new apartment entry (ID 1286, name "벽산").
$ws.Rows("31:31").Insert()
$ws.Range("A31").Value = 1286
$ws.Range("B31").Value = "벽산"

# Match the saved view state: selection on B32, scrolled so row 7 is the
# first visible row.
$null = $ws.Range("B32").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
